$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.085.86'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4644'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2775'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06391'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.17'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +13.87%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07529'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.838.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.970'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.6249'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '293.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +20.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '30.027.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007366'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.084.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.082'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.069'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.936'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1074'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.005'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.812'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04895'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7257'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.111'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.728'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01895'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.653'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.968'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8580'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.64'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.670'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4038'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.046'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.966'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.43%  '
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('E51').Value = '  +0.19%  '
